$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Anthony Edwards -> Collin Sexton
$ws.Range("A3").Value = "Collin Sexton"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Utah Jazz"

# Row 7: Kyle Kuzma -> Giannis Antetokounmpo
$ws.Range("A7").Value = "Giannis Antetokounmpo"
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Milwaukee Bucks"

# Row 14: Collin Sexton -> Payton Pritchard
$ws.Range("A14").Value = "Payton Pritchard"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "Boston Celtics"

# Row 15: Harrison Barnes -> Kyle Kuzma
$ws.Range("A15").Value = "Kyle Kuzma"
$ws.Range("B15").Value = "PF"
$ws.Range("C15").Value = "Washington Wizards"

# Row 16: Giannis Antetokounmpo -> Anthony Edwards
$ws.Range("A16").Value = "Anthony Edwards"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Minnesota Timberwolves"
